# Removendo o validação de valor do recibo no banco de dados
#
# Populate Sheet1 with the recibo (receipt) header/data table:
#   Row 1 - bold, centered, thin-bordered column headers
#   Row 2 - the corresponding values for this receipt

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Header row (A1:F1) ------------------------------------------------
$headers = New-Object 'object[,]' 1,6
$headers[0,0] = "Cliente"
$headers[0,1] = "Profissional"
$headers[0,2] = "Descrição do Serviço"
$headers[0,3] = "Forma de Pagamento"
$headers[0,4] = "Data do Pagamento"
$headers[0,5] = "SubTotal"
$ws.Range("A1:F1").Value = $headers

# -- Data row (A2:F2) ---------------------------------------------------
$values = New-Object 'object[,]' 1,6
$values[0,0] = "Basic System"
$values[0,1] = "Desenvolvedor: Maria Lima"
$values[0,2] = "Raspagem de Dados, Leitura PDF"
$values[0,3] = "Pagamento na Entrega"
$values[0,4] = "30/08/2023"
$values[0,5] = "470,00"
$ws.Range("A2:F2").Value = $values

# -- Style the header row: bold font, thin box border, centered/top ----
$a1 = $ws.Cells.Item(1, 1)
$a1.Font.Bold = $true
$a1.Borders.LineStyle = 1
$a1.HorizontalAlignment = -4108
$a1.VerticalAlignment = -4160

# Propagate the exact same computed style to the rest of the header
# row via a format-only paste so no extra intermediate styles are
# created (mirrors what Excel does with the Format Painter).
$a1.Copy() | Out-Null
$ws.Range("B1:F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
